$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Level 1 Testing (sheet2): add "Priority" column, new rows TC_6 / TC_7
# ---------------------------------------------------------------------
$wsL1 = $wb.Worksheets.Item("Level 1 Testing")

$wsL1.Range("E1").Value = "Priority"
$wsL1.Range("E2").Value = "High"
$wsL1.Range("E4").Value = "Medium"
$wsL1.Range("E6").Value = "Medium"
$wsL1.Range("E8").Value = "High"
$wsL1.Range("E10").Value = "High"

$wsL1.Range("A12").Value = "TC_6"
$wsL1.Range("B12").Value = "Level 1 > Timer"
$wsL1.Range("C12").Value = "After Timer finishes player should lose the game"
$wsL1.Range("D12").Value = "Timer is going in negative cycle."
$wsL1.Range("E12").Value = "Medium"
$wsL1.Range("A12:E12").WrapText = $true
$wsL1.Rows.Item(12).RowHeight = 43.2

$wsL1.Range("A14").Value = "TC_7"
$wsL1.Range("B14").Value = "Level 1 > Finish line"
$wsL1.Range("C14").Value = "After Level 1 finishes, Level 1 Finish message should pop-up"
$wsL1.Range("D14").Value = "No message"
$wsL1.Range("E14").Value = "Low"
$wsL1.Range("A14:E14").WrapText = $true
$wsL1.Rows.Item(14).RowHeight = 43.2

# ---------------------------------------------------------------------
# Level 2 Testing (sheet3): add "Priority" column, new rows TC_3/4/5
# ---------------------------------------------------------------------
$wsL2 = $wb.Worksheets.Item("Level 2 Testing")

$wsL2.Range("E1").Value = "Priority"
$wsL2.Range("E2").Value = "High"
$wsL2.Range("E4").Value = "High"

$wsL2.Range("A6").Value = "TC_3"
$wsL2.Range("B6").Value = "Level 2 > Timer"
$wsL2.Range("C6").Value = "After Timer finishes player should lose the game"
$wsL2.Range("D6").Value = "Timer is going in negative cycle."
$wsL2.Range("E6").Value = "High"
$wsL2.Range("A6:E6").WrapText = $true
$wsL2.Rows.Item(6).RowHeight = 28.8

$wsL2.Range("A8").Value = "TC_4"
$wsL2.Range("B8").Value = "Level 2 > Spikes"
$wsL2.Range("C8").Value = "Spikes should reduce HP"
$wsL2.Range("D8").Value = "HP reduction functionality needs to be implemented"
$wsL2.Range("E8").Value = "High"
$wsL2.Range("A8:E8").WrapText = $true
$wsL2.Rows.Item(8).RowHeight = 28.8

$wsL2.Range("A10").Value = "TC_5"
$wsL2.Range("B10").Value = "Level 2 > Finish line"
$wsL2.Range("C10").Value = "After Level 2 finishes, Level 2 Finish message should pop-up"
$wsL2.Range("D10").Value = "No message"
$wsL2.Range("E10").Value = "Low"
$wsL2.Range("A10:E10").WrapText = $true
$wsL2.Rows.Item(10).RowHeight = 43.2

# ---------------------------------------------------------------------
# Level 3 Testing (sheet4): add "Priority" column, new rows TC_3/7/8
# ---------------------------------------------------------------------
$wsL3 = $wb.Worksheets.Item("Level 3 Testing")

$wsL3.Range("E2").Value = "High"
$wsL3.Range("E4").Value = "High"

$wsL3.Range("A6").Value = "TC_3"
$wsL3.Range("B6").Value = "Level 3 > Timer"
$wsL3.Range("C6").Value = "After Timer finishes player should lose the game"
$wsL3.Range("D6").Value = "Timer is going in negative cycle."
$wsL3.Range("E6").Value = "High"
$wsL3.Range("A6:E6").WrapText = $true
$wsL3.Rows.Item(6).RowHeight = 43.2

$wsL3.Range("A8").Value = "TC_7"
$wsL3.Range("B8").Value = "Level 3 > Finish line"
$wsL3.Range("C8").Value = "After Level 1 finishes, Level 1 Finish message should pop-up"
$wsL3.Range("D8").Value = "No message"
$wsL3.Range("E8").Value = "Low"
$wsL3.Range("A8:E8").WrapText = $true
$wsL3.Rows.Item(8).RowHeight = 43.2

$wsL3.Range("A10").Value = "TC_8"
$wsL3.Range("B10").Value = "Level 3 > Difficulty level"
$wsL3.Range("C10").Value = "More difficulties needed"
$wsL3.Range("D10").Value = "Less Enemy"
$wsL3.Range("E10").Value = "Medium"

# ---------------------------------------------------------------------
# Selections / active cells to match the end-of-edit UI state
# ---------------------------------------------------------------------
$wsL2.Rows.Item(10).Select()
$wsL2.Range("A10:E10").Select()
$wsL2.Range("A10").Activate()

$wsL3.Range("E10").Select()

$wsL1.Rows.Item(14).Select()

$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Activate()
$ws1.Range("C10").Select()
